$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = '65.524.06'
$ws.Range("D3").Value = '2.647.45'
$ws.Range("E3").Value = '  -1.23%  '
$ws.Range("E4").Value = '  +0.05%  '
Set-TextValue "D5" '595.90'
$ws.Range("E5").Value = '  -1.32%  '
Set-TextValue "D6" '155.46'
$ws.Range("E6").Value = '  -0.78%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +6.38%  '
$ws.Range("E9").Value = '  +3.05%  '
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("E11").Value = '  -2.81%  '
$ws.Range("E12").Value = '  +0.25%  '
Set-TextValue "D13" '28.72'
$ws.Range("E13").Value = '  -2.58%  '
$ws.Range("E14").Value = '  -2.93%  '
$ws.Range("D15").Value = '3.123.00'
$ws.Range("E15").Value = '  -1.20%  '
$ws.Range("D16").Value = '65.400.25'
$ws.Range("E16").Value = '  -0.34%  '
$ws.Range("D17").Value = '2.664.90'
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("E18").Value = '  -0.05%  '
Set-TextValue "D19" '4.74'
$ws.Range("E19").Value = '  -1.80%  '
$ws.Range("E20").Value = '  -1.83%  '
Set-TextValue "D21" '348.08'
$ws.Range("E21").Value = '  -1.17%  '
Set-TextValue "D22" '1.00'
$ws.Range("E22").Value = '  +0.08%  '
Set-TextValue "D23" '68.97'
$ws.Range("E23").Value = '  -2.08%  '
$ws.Range("E24").Value = '  +1.20%  '
Set-TextValue "D25" '9.59'
$ws.Range("E25").Value = '  -2.59%  '
$ws.Range("E26").Value = '  +0.57%  '
$ws.Range("E27").Value = '  -3.14%  '
$ws.Range("E28").Value = '  -3.24%  '
$ws.Range("E29").Value = '  +0.09%  '
Set-TextValue "D30" '537.58'
$ws.Range("E30").Value = '  +0.19%  '
Set-TextValue "D31" '7.90'
$ws.Range("E31").Value = '  -3.36%  '
Set-TextValue "D32" '2.13'
$ws.Range("E32").Value = '  -1.58%  '
$ws.Range("E33").Value = '  -1.10%  '
$ws.Range("E34").Value = '  -2.84%  '
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("E36").Value = '  -1.79%  '
$ws.Range("E37").Value = '  -1.15%  '
$ws.Range("E38").Value = '  +0.08%  '
Set-TextValue "D39" '155.09'
$ws.Range("E39").Value = '  -3.38%  '
$ws.Range("E40").Value = '  -3.12%  '
$ws.Range("E41").Value = '  -0.01%  '
Set-TextValue "D42" '160.53'
$ws.Range("E42").Value = '  -3.53%  '
$ws.Range("E44").Value = '  +2.97%  '
Set-TextValue "D45" '0.0602'
$ws.Range("E45").Value = '  -2.93%  '
Set-TextValue "D46" '22.44'
$ws.Range("E46").Value = '  -3.04%  '
$ws.Range("E47").Value = '  -2.42%  '
$ws.Range("E48").Value = '  -3.72%  '
Set-TextValue "D49" '0.100'
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("D50").Value = '0.0₆0252'
$ws.Range("E50").Value = '  +7.34%  '
Set-TextValue "D51" '19.58'
$ws.Range("E51").Value = '  -3.77%  '
